# Scheduled market-data refresh: updates currentAveragePrice / LevePrice /
# LeveProfit columns (H-N) for the affected Leve rows across each job sheet,
# reflecting latest market board pulls. A few rows gain/lose a profit cell
# (HQ vs NQ) as the cheaper recipe side flips.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Range("H137").Value = 7027.0625
$ws.Range("I137").Value = 7082.273
$ws.Range("J137").Value = 6905.6
$ws.Range("K137").Value = 21246.819
$ws.Range("L137").Value = 20716.8
$ws.Range("M137").Value = -18696.819
$ws.Range("N137").Value = -25816.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2192.1177
$ws.Range("I32").Value = 2192.1177
$ws.Range("K32").Value = 2192.1177
$ws.Range("M32").Value = -1905.1177
# Row 61
$ws.Range("H61").Value = 3309.4614
$ws.Range("I61").Value = 3090.1
$ws.Range("J61").Value = 4040.6667
$ws.Range("K61").Value = 3090.1
$ws.Range("L61").Value = 4040.6667
$ws.Range("M61").Value = -2878.1
$ws.Range("N61").Value = -4464.6667
# Row 63
$ws.Range("H63").Value = 51500
$ws.Range("I63").Value = 51500
$ws.Range("K63").Value = 51500
$ws.Range("M63").Value = -50814
# Row 66
$ws.Range("H66").Value = 51500
$ws.Range("I66").Value = 51500
$ws.Range("K66").Value = 257500
$ws.Range("M66").Value = -254068
# Row 74
$ws.Range("H74").Value = 1099.7273
$ws.Range("I74").Value = 1099.7273
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1099.7273
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -225.7273
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 1099.7273
$ws.Range("I77").Value = 1099.7273
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5498.636500000001
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1130.636500000001
$ws.Range("N77").ClearContents()
# Row 97
$ws.Range("H97").Value = 1741.4
$ws.Range("I97").Value = 823.7778
$ws.Range("J97").Value = 10000
$ws.Range("K97").Value = 823.7778
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = -327.7778
$ws.Range("N97").Value = -10992
# Row 110
$ws.Range("H110").Value = 6043.696
$ws.Range("I110").Value = 5945
$ws.Range("J110").Value = 6399
$ws.Range("K110").Value = 5945
$ws.Range("L110").Value = 6399
$ws.Range("M110").Value = -3900
$ws.Range("N110").Value = -10489
# Row 125
$ws.Range("H125").Value = 94166.5
$ws.Range("J125").Value = 94166.5
$ws.Range("L125").Value = 94166.5
$ws.Range("N125").Value = -104006.5
# Row 136
$ws.Range("H136").Value = 3309.4614
$ws.Range("I136").Value = 3090.1
$ws.Range("J136").Value = 4040.6667
$ws.Range("K136").Value = 9270.299999999999
$ws.Range("L136").Value = 12122.0001
$ws.Range("M136").Value = -6720.299999999999
$ws.Range("N136").Value = -17222.0001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1923.5883
$ws.Range("I20").Value = 1927.1818
$ws.Range("K20").Value = 1927.1818
$ws.Range("M20").Value = -1680.1818
# Row 80
$ws.Range("H80").Value = 6140.1665
$ws.Range("I80").Value = 95.8
$ws.Range("J80").Value = 10457.571
$ws.Range("K80").Value = 95.8
$ws.Range("L80").Value = 10457.571
$ws.Range("M80").Value = 902.2
$ws.Range("N80").Value = -12453.571
# Row 83
$ws.Range("H83").Value = 6140.1665
$ws.Range("I83").Value = 95.8
$ws.Range("J83").Value = 10457.571
$ws.Range("K83").Value = 479
$ws.Range("L83").Value = 52287.855
$ws.Range("M83").Value = 4513
$ws.Range("N83").Value = -62271.855
# Row 94
$ws.Range("H94").Value = 890
$ws.Range("I94").Value = 953.3333
$ws.Range("J94").Value = 805.55554
$ws.Range("K94").Value = 953.3333
$ws.Range("L94").Value = 805.55554
$ws.Range("M94").Value = -502.3333
$ws.Range("N94").Value = -1707.55554

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2084.6365
$ws.Range("J16").Value = 2789.2856
$ws.Range("L16").Value = 2789.2856
$ws.Range("N16").Value = -3363.2856
# Row 31
$ws.Range("H31").Value = 2089.5
$ws.Range("I31").Value = 2137.6
$ws.Range("J31").Value = 1849
$ws.Range("K31").Value = 2137.6
$ws.Range("L31").Value = 1849
$ws.Range("M31").Value = -1842.6
$ws.Range("N31").Value = -2439
# Row 34
$ws.Range("H34").Value = 2089.5
$ws.Range("I34").Value = 2137.6
$ws.Range("J34").Value = 1849
$ws.Range("K34").Value = 2137.6
$ws.Range("L34").Value = 1849
$ws.Range("M34").Value = -1935.6
$ws.Range("N34").Value = -2253
# Row 96
$ws.Range("H96").Value = 11835
$ws.Range("J96").Value = 11835
$ws.Range("L96").Value = 11835
$ws.Range("N96").Value = -17327
# Row 113
$ws.Range("H113").Value = 2084.6365
$ws.Range("J113").Value = 2789.2856
$ws.Range("L113").Value = 2789.2856
$ws.Range("N113").Value = -7129.2856
# Row 141
$ws.Range("H141").Value = 169999
$ws.Range("J141").Value = 169999
$ws.Range("L141").Value = 169999
$ws.Range("N141").Value = -180359

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 140
$ws.Range("H140").Value = 627181.5
$ws.Range("I140").Value = 627181.5
$ws.Range("K140").Value = 1881544.5
$ws.Range("M140").Value = -1876364.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2899.5
$ws.Range("I80").Value = 2899.5
$ws.Range("K80").Value = 2899.5
$ws.Range("M80").Value = -1901.5
# Row 83
$ws.Range("H83").Value = 2899.5
$ws.Range("I83").Value = 2899.5
$ws.Range("K83").Value = 14497.5
$ws.Range("M83").Value = -9505.5
# Row 97
$ws.Range("H97").Value = 698.4
$ws.Range("I97").Value = 698.4
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 698.4
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -202.4
$ws.Range("N97").ClearContents()
# Row 132
$ws.Range("H132").Value = 2267.077
$ws.Range("I132").Value = 2297.5454
$ws.Range("K132").Value = 6892.6362
$ws.Range("M132").Value = -4362.6362

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1548.5
$ws.Range("I22").Value = 1647.25
$ws.Range("J22").Value = 1351
$ws.Range("K22").Value = 1647.25
$ws.Range("L22").Value = 1351
$ws.Range("M22").Value = -1352.25
$ws.Range("N22").Value = -1941
# Row 27
$ws.Range("H27").Value = 1548.5
$ws.Range("I27").Value = 1647.25
$ws.Range("J27").Value = 1351
$ws.Range("K27").Value = 1647.25
$ws.Range("L27").Value = 1351
$ws.Range("M27").Value = -1540.25
$ws.Range("N27").Value = -1565
# Row 61
$ws.Range("H61").Value = 1383.8
$ws.Range("I61").Value = 1424.5
$ws.Range("K61").Value = 1424.5
$ws.Range("M61").Value = -1222.5
# Row 113
$ws.Range("H113").Value = 1383.8
$ws.Range("I113").Value = 1424.5
$ws.Range("K113").Value = 1424.5
$ws.Range("M113").Value = 745.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 139990
$ws.Range("I29").Value = 139990
$ws.Range("K29").Value = 139990
$ws.Range("M29").Value = -139700
# Row 95
$ws.Range("H95").Value = 24172
$ws.Range("J95").Value = 24172
$ws.Range("L95").Value = 24172
$ws.Range("N95").Value = -29664
# Row 136
$ws.Range("H136").Value = 8166.524
$ws.Range("I136").Value = 8166.524
$ws.Range("K136").Value = 24499.572
$ws.Range("M136").Value = -21949.572
